# Generate Report for Handback
# - Update the "Status" text (shared across Overview zh-cn/de-de columns and
#   the per-locale "Status" table column) from "In Translation" to
#   "Handed back: in sync with en-US".
# - Fill in the "Latest Target File" / "Latest Handback File" /
#   "Latest Handback DateTime" columns for both locale sheets (previously
#   blank / placeholder "0001-01-01 00:00:00"), and hyperlink the target
#   file name to the same source-doc URL as column A.
# - Widen the now-longer "Status" / "Latest Target File" / "Latest Handback
#   File" columns so the new text is not clipped.

$wb = $excel.ActiveWorkbook

$statusOld = "In Translation"
$statusNew = "Handed back: in sync with en-US"

$mdUrlB7 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/84692ea7b98921da544303388af8612021dc800a/e2e/b7cc3251-cc38-40d7-a6a2-d425982365d5.md"
$mdUrlD8 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/84692ea7b98921da544303388af8612021dc800a/e2e/d86d0e51-b6b8-4c60-af8e-1592752585e5.md"

$mdNameB7 = "b7cc3251-cc38-40d7-a6a2-d425982365d5.md"
$mdNameD8 = "d86d0e51-b6b8-4c60-af8e-1592752585e5.md"

# ---------------------------------------------------------------------
# Overview sheet: refresh the "zh-cn" / "de-de" status cells.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew

# AutoFit, then nudge to the known post-AutoFit width for this text so the
# stored column width matches what Excel computes for the longer string.
$wsOverview.Columns.Item(5).AutoFit()
$wsOverview.Columns.Item(6).AutoFit()
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet (row 2 = b7cc3251..., row 3 = d86d0e51...)
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew

# Row 2 - Latest Target File / Latest Handback File / Latest Handback DateTime
$wsZh.Range("I2").Value = $mdNameB7
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrlB7, "", "", $mdNameB7) | Out-Null
$wsZh.Range("J2").Value = "b7cc3251-cc38-40d7-a6a2-d425982365d5.01daa9a24eaa4aac00812c0d0d65af58d246ee7b.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-03 12:27:22"

# Row 3
$wsZh.Range("I3").Value = $mdNameD8
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrlD8, "", "", $mdNameD8) | Out-Null
$wsZh.Range("J3").Value = "d86d0e51-b6b8-4c60-af8e-1592752585e5.3e36770600ebae4da649200134152ced03bd12ca.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-03 12:27:22"

$wsZh.Columns.Item(3).AutoFit()
$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet (row 2 = b7cc3251..., row 3 = d86d0e51...)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

# Row 2
$wsDe.Range("I2").Value = $mdNameB7
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrlB7, "", "", $mdNameB7) | Out-Null
$wsDe.Range("J2").Value = "b7cc3251-cc38-40d7-a6a2-d425982365d5.01daa9a24eaa4aac00812c0d0d65af58d246ee7b.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-03 12:27:28"

# Row 3
$wsDe.Range("I3").Value = $mdNameD8
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrlD8, "", "", $mdNameD8) | Out-Null
$wsDe.Range("J3").Value = "d86d0e51-b6b8-4c60-af8e-1592752585e5.3e36770600ebae4da649200134152ced03bd12ca.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-03 12:27:28"

$wsDe.Columns.Item(3).AutoFit()
$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Host "Handback report generated."
